# Applies the "Handles float input without breaking stuff" edit to the
# marksheet workbook: fills in the student's answers (previously blank /
# "Absent"), recomputes the summary block, and drops the now-unused third
# answer-block columns (G:H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Summary block (rows 10-12): Right / Wrong / Not-Attempt / Max, and
#    the marking scheme + totals. Row labels (A10:A12) pick up the same
#    "mtitleStyle" heading style as A9.
# ---------------------------------------------------------------------
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A10:A12").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B10").Value2 = 20
$ws.Range("C10").Value2 = 3
$ws.Range("D10").Value2 = 5
$ws.Range("E10").Value2 = 28

$ws.Range("B11").Value2 = 4
$ws.Range("C11").Value2 = -1

$ws.Range("B12").Value2 = 80
$ws.Range("C12").Value2 = -3
$ws.Range("E12").Value2 = "77/112"

# ---------------------------------------------------------------------
# 2. Per-question student answers (column A, rows 16-40). Style mirrors
#    correctness: blank -> normalStyle, match with column B ->
#    correctStyle, mismatch -> incorrectStyle.
# ---------------------------------------------------------------------
$answers = @{
  16 = "Option A"
  17 = $null
  18 = "Option B"
  19 = "Option C"
  20 = "Option B"
  21 = "Option C"
  22 = "Option D"
  23 = $null
  24 = $null
  25 = "Option A"
  26 = "Option C"
  27 = "Option A"
  28 = "Option D"
  29 = $null
  30 = "Option B"
  31 = "Option D"
  32 = "Option C"
  33 = "Option D"
  34 = "Option A"
  35 = $null
  36 = "Option A"
  37 = "Option B"
  38 = "Option A"
  39 = "Option D"
  40 = "Option B"
}

$ws.Range("B10").Copy() | Out-Null   # correctStyle source
foreach ($r in @(16,18,19,20,21,22,25,26,27,28,30,31,32,33,36,38,39)) {
  $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
}

$ws.Range("C10").Copy() | Out-Null   # incorrectStyle source
foreach ($r in @(34,37,40)) {
  $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
}
# rows 17,23,24,29,35 keep their existing blank / normalStyle formatting.

foreach ($r in 16..40) {
  $v = $answers[$r]
  if ($v -ne $null) {
    $ws.Range("A$r").Value2 = $v
  }
}

# ---------------------------------------------------------------------
# 3. Row 16-18 "D" column picks up the corresponding "E" (correct
#    answer) value with the correctStyle formatting, matching column A.
# ---------------------------------------------------------------------
$ws.Range("B10").Copy() | Out-Null   # correctStyle source
$ws.Range("D16:D18").PasteSpecial(-4122) | Out-Null
$ws.Range("D16").Value2 = $ws.Range("E16").Value2
$ws.Range("D17").Value2 = $ws.Range("E17").Value2
$ws.Range("D18").Value2 = $ws.Range("E18").Value2

# ---------------------------------------------------------------------
# 4. Drop the unused third answer-block columns (G:H) entirely, which
#    also removes D/E for rows 19-40 (only rows 15-18 keep them).
# ---------------------------------------------------------------------
$ws.Range("D19:E40").Clear() | Out-Null
$ws.Columns("G:H").Delete() | Out-Null

$excel.CutCopyMode = 0
